$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Ascelline" row (row 114: Name=Ascelline, Ch.=20, ID=0x56, Class & Level=???,
# Method=?) is removed entirely. Excel's EntireRow delete shifts every row below it
# up by one, which reproduces the rest of the diff (rows 116-167 becoming 115-166)
# without needing to touch any other cell individually.
$ws.Rows(114).Delete()

# Update the view state to match the author's saved selection/scroll position.
$excel.ActiveWindow.ScrollRow = 145
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A158").Select()
